$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its values as text so Excel does not
# reinterpret numeric-looking strings (e.g. "314.84", "0.00001343") as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$data = @(
    @{Row=2; B='Bitcoin'; C='https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'; D='24.862.79'; E='  +1.04%  '}
    @{Row=3; B='Ethereum'; C='https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'; D='1.708.52'; E='  +1.19%  '}
    @{Row=4; B='TetherUSD'; C='https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'; D='1.003'; E='  -0.25%  '}
    @{Row=5; B='BNB'; C='https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'; D='314.84'; E='  +0.14%  '}
    @{Row=6; B='USDC'; C='https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'; D='1.003'; E='  -0.44%  '}
    @{Row=7; B='XRP'; C='https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'; D='0.4022'; E='  +3.11%  '}
    @{Row=8; B='Cardano'; C='https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'; D='0.4050'; E='  +0.67%  '}
    @{Row=9; B='BinanceUSD'; C='https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D='1.002'; E='  -0.47%  '}
    @{Row=10; B='OKB'; C='https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; D='53.72'; E='  +1.76%  '}
    @{Row=11; B='Polygon'; C='https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; D='1.472'; E='  -0.90%  '}
    @{Row=12; B='Dogecoin'; C='https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; D='0.08801'; E='  +1.09%  '}
    @{Row=13; B='Solana'; C='https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D='26.28'; E='  +7.13%  '}
    @{Row=14; B='Polkadot'; C='https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D='7.510'; E='  -1.49%  '}
    @{Row=15; B='Chainlink'; C='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D='8.004'; E='  +0.75%  '}
    @{Row=16; B='ShibaInu'; C='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D='0.00001343'; E='  +0.70%  '}
    @{Row=17; B='WrappedEther'; C='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D='1.711.35'; E='  +1.44%  '}
    @{Row=18; B='Litecoin'; C='https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D='95.52'; E='  -2.49%  '}
    @{Row=19; B='TRON'; C='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D='0.07177'; E='  +1.19%  '}
    @{Row=20; B='Avalanche'; C='https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D='20.97'; E='  +6.84%  '}
    @{Row=21; B='Uniswap'; C='https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D='7.292'; E='  +0.42%  '}
    @{Row=22; B='Dai'; C='https://coinranking.com/coin/MoTuySvg7+dai-dai'; D='1.004'; E='  -0.22%  '}
    @{Row=23; B='Cosmos'; C='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D='14.47'; E='  +2.21%  '}
    @{Row=24; B='WrappedBTC'; C='https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D='24.863.56'; E='  +1.10%  '}
    @{Row=25; B='Toncoin'; C='https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D='2.337'; E='  -0.55%  '}
    @{Row=26; B='LidoDAOToken'; C='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D='2.886'; E='  -3.87%  '}
    @{Row=27; B='HuobiToken'; C='https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; D='6.413'; E='  +22.72%  '}
    @{Row=28; B='EthereumClassic'; C='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D='23.07'; E='  +2.24%  '}
    @{Row=29; B='Monero'; C='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D='161.71'; E='  +0.20%  '}
    @{Row=30; B='BitcoinCash'; C='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D='144.12'; E='  +6.00%  '}
    @{Row=31; B='Filecoin'; C='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D='8.379'; E='  -1.58%  '}
    @{Row=32; B='WEMIXTOKEN'; C='https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; D='2.281'; E='  +15.44%  '}
    @{Row=33; B='WrappedliquidstakedEther2.0'; C='https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; D='1.870.24'; E='  +0.04%  '}
    @{Row=34; B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.08727'; E='  -0.43%  '}
    @{Row=35; B='VeChain'; C='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D='0.03184'; E='  +10.03%  '}
    @{Row=36; B='InternetComputer(DFINITY)'; C='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D='7.187'; E='  -3.63%  '}
    @{Row=37; B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='1.028'; E='  -0.60%  '}
    @{Row=38; B='Algorand'; C='https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; D='0.2870'; E='  +6.07%  '}
    @{Row=39; B='TheSandbox'; C='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; D='0.8422'; E='  +9.02%  '}
    @{Row=40; B='FraxShare'; C='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D='10.83'; E='  +1.24%  '}
    @{Row=41; B='Stellar'; C='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D='0.09455'; E='  +3.99%  '}
    @{Row=42; B='Aptos'; C='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D='14.20'; E='  +1.37%  '}
    @{Row=43; B='TrustWalletToken'; C='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D='1.482'; E='  +1.98%  '}
    @{Row=44; B='EnergySwap'; C='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D='17.55'; E='  +5.71%  '}
    @{Row=45; B='NEARProtocol'; C='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D='2.719'; E='  +6.15%  '}
    @{Row=46; B='Decentraland'; C='https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'; D='0.7435'; E='  +4.40%  '}
    @{Row=47; B='PancakeSwap'; C='https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; D='4.229'; E='  +0.61%  '}
    @{Row=48; B='Flow'; C='https://coinranking.com/coin/QQ0NCmjVq+flow-flow'; D='1.376'; E='  +2.82%  '}
    @{Row=49; B='Frax'; C='https://coinranking.com/coin/KfWtaeV1W+frax-frax'; D='1.001'; E='  -0.51%  '}
    @{Row=50; B='Quant'; C='https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; D='140.59'; E='  +1.71%  '}
    @{Row=51; B='Cronos'; C='https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; D='0.08392'; E='  +5.75%  '}
)

foreach ($row in $data) {
    $ws.Cells.Item($row.Row, 2).Value = $row.B
    $ws.Cells.Item($row.Row, 3).Value = $row.C
    $ws.Cells.Item($row.Row, 4).Value = $row.D
    $ws.Cells.Item($row.Row, 5).Value = $row.E
}

$wb.Save()
